# Update (Analyze PO & Forecast)
# The forecast was re-run 4 weeks earlier than before: every Week_Start_Date
# on the "Forecast Comparison" sheet shifts back 28 days, the first two
# MyForecast values change to reflect the new run, the is_holiday_week flags
# are cleared (no longer computed), and the "Summary" sheet's derived totals
# are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: Week_Start_Date (col B), shifted back 28 days ---
$newDates = @{
    2  = "2024-12-29"
    3  = "2025-01-05"
    4  = "2025-01-12"
    5  = "2025-01-19"
    6  = "2025-01-26"
    7  = "2025-02-02"
    8  = "2025-02-09"
    9  = "2025-02-16"
    10 = "2025-02-23"
    11 = "2025-03-02"
    12 = "2025-03-09"
    13 = "2025-03-16"
    14 = "2025-03-23"
    15 = "2025-03-30"
    16 = "2025-04-06"
    17 = "2025-04-13"
}

foreach ($row in $newDates.Keys) {
    # Leading apostrophe keeps this a literal text value (matches the
    # original inline-string storage) instead of being parsed as a date.
    $ws.Range("B$row").Value = "'" + $newDates[$row]
}

# --- Forecast Comparison: MyForecast (col D) for the first two weeks ---
$ws.Range("D2").Value = 367
$ws.Range("D3").Value = 356

# --- Forecast Comparison: is_holiday_week (col J) no longer populated ---
$ws.Range("J2:J17").ClearContents()

# --- Summary: refreshed aggregate figures ---
$summary.Range("B9").Value  = "'4012"          # Total Forecast (16 Weeks)
$summary.Range("B10").Value = "'2218"          # Total Forecast (8 Weeks)
$summary.Range("B11").Value = "'1238"          # Total Forecast (4 Weeks)
$summary.Range("B12").Value = "'367"           # Max Forecast
$summary.Range("B13").Value = "'2024-12-29"    # Max Forecast Week
$summary.Range("B15").Value = "'2025-04-13"    # Min Forecast Week
